{"js": "// Insert a new paragraph \"Given a set of planes with target landing times and \"\n// right after the \"Problem Context\" Heading 1 paragraph (before the blank\n// paragraph that already follows it), matching the new paragraph's plain\n// (Normal) formatting \u2014 only `lang=\"en-US\"` on the run/paragraph marks.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\n// Find the \"Problem Context\" heading paragraph.\nlet headingIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text.trim() === \"Problem Context\" && /Heading 1/i.test(para.style)) {\n    headingIndex = i;\n    break;\n  }\n}\n\nif (headingIndex === -1) {\n  throw new Error('Could not find the \"Problem Context\" heading paragraph.');\n}\n\nconst newText = \"Given a set of planes with target landing times and \";\n\n// Insert the new paragraph so it lands right after the heading. Inserting\n// \"Before\" the paragraph that currently follows the heading (rather than\n// \"After\" the heading itself) means the new paragraph inherits that\n// following paragraph's plain formatting instead of the Heading 1 style.\nlet newParagraph;\nif (headingIndex + 1 < paragraphs.items.length) {\n  const following = paragraphs.items[headingIndex + 1];\n  newParagraph = following.insertParagraph(newText, Word.InsertLocation.before);\n} else {\n  // Fallback: no following paragraph, just insert after the heading.\n  newParagraph = paragraphs.items[headingIndex].insertParagraph(newText, Word.InsertLocation.after);\n}\n\nawait context.sync();\n", "ps1": "# Insert a new paragraph \"Given a set of planes with target landing times and \"\n# right after the \"Problem Context\" Heading 1 paragraph (before the blank\n# paragraph that already follows it), so the new paragraph ends up with\n# plain (Normal) formatting instead of inheriting the Heading 1 style.\n\n$d = $word.ActiveDocument\n\n# Locate the \"Problem Context\" heading (case/whole-word sensitive so we\n# don't match the lower-case \"problem context\" mention earlier in the text).\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.MatchCase = $true\n$rng.Find.MatchWholeWord = $true\n$found = $rng.Find.Execute(\"Problem Context\")\nif (-not $found) {\n    throw \"Could not find the 'Problem Context' heading paragraph.\"\n}\n\n# Resolve which paragraph (1-based index into $d.Paragraphs) the found\n# range falls in, then confirm it's really the Heading 1 occurrence.\n$headingIndex = 0\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($rng.Start -ge $p.Range.Start -and $rng.Start -lt $p.Range.End) {\n        $headingIndex = $i\n        break\n    }\n}\nif ($headingIndex -eq 0 -or $d.Paragraphs.Item($headingIndex).Style.NameLocal -ne \"Heading 1\") {\n    throw \"Could not resolve the 'Problem Context' heading paragraph.\"\n}\n\n$newText = \"Given a set of planes with target landing times and \"\n\nif ($headingIndex -lt $d.Paragraphs.Count) {\n    # Insert before the paragraph that currently follows the heading so the\n    # new paragraph inherits that (plain) paragraph's formatting.\n    $followingPara = $d.Paragraphs.Item($headingIndex + 1)\n    $followingPara.Range.InsertParagraphBefore()\n    $insertedPara = $d.Paragraphs.Item($headingIndex + 1)\n    $insertedPara.Range.Text = $newText\n} else {\n    # Fallback: heading is the last paragraph in the document.\n    $headingPara = $d.Paragraphs.Item($headingIndex)\n    $headingPara.Range.InsertParagraphAfter()\n    $insertedPara = $d.Paragraphs.Item($headingIndex + 1)\n    $insertedPara.Range.Text = $newText\n}\n"}
